$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Frequency (column B) and Probability (column C) values for rows 2-11
$values = @(
    @(21, 0.105),
    @(13, 0.065),
    @(20, 0.1),
    @(21, 0.105),
    @(20, 0.1),
    @(25, 0.125),
    @(15, 0.075),
    @(24, 0.12),
    @(20, 0.1),
    @(21, 0.105)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}
